$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 5.834599999999997
$ws.Range("E4").Value = 13.59290000000001
$ws.Range("E5").Value = 12.9632
$ws.Range("B7").Value = 6.765399999999996
$ws.Range("E8").Value = 14.0748
$ws.Range("B16").Value = 8.975600000000009
$ws.Range("E16").Value = 12.56310000000001

$wb.Save()
